$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 125: new result columns + updated odds ----
$ws.Range("H125").Value = 2
$ws.Range("I125").Value = 1
$ws.Range("J125").Value = "H"

$ws.Range("N125").Value = 1.8
$ws.Range("O125").Value = 3.75
$ws.Range("P125").Value = 4.2
$ws.Range("Q125").Value = -0.75
$ws.Range("R125").Value = 2.025
$ws.Range("S125").Value = 1.825
$ws.Range("U125").Value = 1.825
$ws.Range("V125").Value = 2.025
$ws.Range("W125").Value = 0.8
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = -1
$ws.Range("Z125").Value = 0.5125
$ws.Range("AA125").Value = -0.5
$ws.Range("AB125").Value = 0.4125
$ws.Range("AC125").Value = -0.5

# ---- Row 126: updated odds ----
$ws.Range("N126").Value = 4.5
$ws.Range("O126").Value = 3.75
$ws.Range("P126").Value = 1.727
$ws.Range("T126").Value = 2.75
$ws.Range("U126").Value = 1.85
$ws.Range("V126").Value = 1.95

# ---- Row 127: updated odds ----
$ws.Range("O127").Value = 3.5
$ws.Range("P127").Value = 3.1
$ws.Range("R127").Value = 2
$ws.Range("S127").Value = 1.8

# ---- Row 128: brand-new match row ----
# Seed formatting by cloning row 127's formats first (keeps A128 bold/
# bordered/centered like the other "id" cells, and E128 carrying the
# date/time number format), then fill in the values.
$ws.Range("A127:G127").Copy()
$ws.Range("A128").PasteSpecial(-4122)

$ws.Range("A128").Value = 126
$ws.Range("B128").Value = 7751765
$ws.Range("C128").Value = "India Super League"
$ws.Range("D128").Value = "India Super League"
$ws.Range("E128").Value = 45392.45833333334
$ws.Range("F128").Value = "Punjab FC"
$ws.Range("G128").Value = "East Bengal Club"

$ws.Range("K128").Value = 2.625
$ws.Range("L128").Value = 3.3
$ws.Range("M128").Value = 2.55
$ws.Range("N128").Value = 2.9
$ws.Range("O128").Value = 3.3
$ws.Range("P128").Value = 2.3
$ws.Range("Q128").Value = 0
$ws.Range("R128").Value = 2.1
$ws.Range("S128").Value = 1.7
$ws.Range("T128").Value = 2.5
$ws.Range("U128").Value = 1.85
$ws.Range("V128").Value = 1.95
$ws.Range("W128").Value = 0
$ws.Range("X128").Value = 0
$ws.Range("Y128").Value = 0
$ws.Range("Z128").Value = 0
$ws.Range("AA128").Value = 0
